$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Prerequisite column (E) with grade qualifiers ("CourseCode:Grade").
# Write order matches the original authoring session so the shared-string
# table is (re)built with the same index assignment as the source edit.
$ws.Range("E3").Value  = "CS112:C"
$ws.Range("E4").Value  = "CS151:C;CS112:C"
$ws.Range("E5").Value  = "CS122:C;CS150:C"
$ws.Range("E6").Value  = "CS122:C"
$ws.Range("E7").Value  = "CS122:C/CS150:C;CS151:C"
$ws.Range("E8").Value  = "CS210:C;CS251:C;MATH241:C"
$ws.Range("E9").Value  = "CS122:C;CS151:C;MATH241:C"
$ws.Range("E10").Value = "CS213:C;MATH241:C;MATH259:C"
$ws.Range("E12").Value = "CS213:C"
$ws.Range("E13").Value = "CS251:C"
$ws.Range("E11").Value = "CS150:C"
$ws.Range("E14").Value = "CS213:C"
$ws.Range("E15").Value = "CS122:C;CS151:C;MATH241:C;CS259:C"
$ws.Range("E16").Value = "CS112:D/CS122:D"
$ws.Range("E17").Value = "CS122:C;CS151:C"
$ws.Range("E21").Value = "MATH122:D"
$ws.Range("E22").Value = "MATH259:C;MATH241:C;CS251:C"
$ws.Range("E24").Value = "ENG101:D"

# Widen column E to fit the longer prerequisite strings (target stored width 43.1640625)
$ws.Columns.Item(5).ColumnWidth = 42.42857142857143

# Restore the view to the originally-recorded scroll position/selection
$excel.Goto($ws.Range("A16"), $true)
$ws.Range("E24").Select()
